{"js": "// Picture Book Final Draft.docx - apply commit changes:\n//  1. Update the due date from \"Friday, May 1\" to \"Wednesday, May 13\"\n//  2. Remove the \"DRAFT - This document is under development\" warning paragraph\n//  3. Remove the decorative horizontal-rule paragraph near the end of the document\n//  4. Resize the two rubric/workshop tables (pct width + narrower columns)\n\n// 1. Update the due date, keeping the existing bold run formatting intact.\nconst dateResults = context.document.body.search(\"Due: Friday, May 1 at 11:59pm\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"Due: Wednesday, May 13 at 11:59pm\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2. Remove the whole \"DRAFT\" warning paragraph (emoji + bold DRAFT + message).\nconst paragraphsForDraft = context.document.body.paragraphs;\nparagraphsForDraft.load(\"items/text\");\nawait context.sync();\n\nlet draftParagraph = null;\nfor (const p of paragraphsForDraft.items) {\n  if (p.text.indexOf(\"DRAFT\") !== -1 && p.text.indexOf(\"under development\") !== -1) {\n    draftParagraph = p;\n    break;\n  }\n}\nif (draftParagraph) {\n  draftParagraph.delete();\n  await context.sync();\n}\n\n// 3. Remove the empty paragraph holding the decorative horizontal-rule VML shape\n//    that sits right after the \"Communicate with your partner\" tip.\nconst paragraphsForHr = context.document.body.paragraphs;\nparagraphsForHr.load(\"items/text\");\nawait context.sync();\n\nlet communicateParagraph = null;\nfor (const p of paragraphsForHr.items) {\n  if (p.text.indexOf(\"Make sure both partners are aligned on the vision and workload distribution\") !== -1) {\n    communicateParagraph = p;\n    break;\n  }\n}\nif (communicateParagraph) {\n  const hrParagraph = communicateParagraph.getNext();\n  hrParagraph.delete();\n  await context.sync();\n}\n\n// 4. Resize both tables: tblW becomes a percentage (4865) and each of the two\n//    columns narrows from 3960 to 3852 twips. There is no supported Table API\n//    member in this runtime that actually persists width changes, so the raw\n//    table OOXML is patched directly via getOoxml()/insertOoxml(replace).\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < tables.items.length; i++) {\n  const table = tables.items[i];\n  const ooxmlResult = table.getOoxml();\n  await context.sync();\n\n  let tableXml = ooxmlResult.value;\n  tableXml = tableXml.replace(\n    '<w:tblW w:w=\"0\" w:type=\"auto\"/>',\n    '<w:tblW w:type=\"pct\" w:w=\"4865\"/>'\n  );\n  tableXml = tableXml.replace(\n    '<w:gridCol w:w=\"3960\"/><w:gridCol w:w=\"3960\"/>',\n    '<w:gridCol w:w=\"3852\"/><w:gridCol w:w=\"3852\"/>'\n  );\n\n  table.insertOoxml(tableXml, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Picture Book Final Draft.docx - apply commit changes:\n#  1. Update the due date from \"Friday, May 1\" to \"Wednesday, May 13\"\n#  2. Remove the \"DRAFT - This document is under development\" warning paragraph\n#  3. Remove the decorative horizontal-rule paragraph near the end of the document\n#  4. Resize the two rubric/workshop tables (pct width + narrower columns)\n#\n# NOTE: table-size edits are done LAST because touching a table's width /\n# column-width properties invalidates paragraph navigation/indexing for the\n# remainder of the script in this runtime, so all paragraph find/delete work\n# must happen first.\n\n$d = $word.ActiveDocument\n\n# 1. Update the due date, keeping the existing bold run formatting intact.\n$find = $d.Content.Find\n$find.Text = \"Due: Friday, May 1 at 11:59pm\"\n$find.Execute()\nif ($find.Found) {\n    $find.Parent.Text = \"Due: Wednesday, May 13 at 11:59pm\"\n}\n\n# 2. Remove the whole \"DRAFT\" warning paragraph (emoji + bold DRAFT + message).\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*DRAFT*\" -and $t -like \"*under development*\") {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 3. Remove the empty paragraph holding the decorative horizontal-rule VML shape\n#    that sits right after the \"Communicate with your partner\" tip.\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Make sure both partners are aligned on the vision and workload distribution*\") {\n        $hrPara = $p.Next()\n        $hrPara.Range.Delete()\n        break\n    }\n}\n\n# 4. Resize both tables: preferred width becomes a percentage (4865/20 = 243.25),\n#    and each of the two columns narrows from 3960 to 3852 twips (3852/20 = 192.6).\nfor ($i = 1; $i -le $d.Tables.Count; $i++) {\n    $t = $d.Tables.Item($i)\n    $t.PreferredWidthType = 2\n    $t.PreferredWidth = 243.25\n    $t.Columns.Item(1).Width = 192.6\n    $t.Columns.Item(2).Width = 192.6\n}\n"}
